# Adiciona circuito quadro elétrico
# Duplicates the "Cronograma" schedule table into a second table
# ("Cronograma da Implementação") and fills in a few missing "X" marks
# on the original ("Cronograma do TCC") table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Duplicate the whole table (title + header + 8 rows) down to row 15 ---
# Do this BEFORE editing table 1 so the new table starts from the same
# baseline formatting the original table had.
$src = $ws.Range("F3:J12")
$dst = $ws.Range("F15")
$src.Copy($dst)

# --- 2. Rename the first (original) table's title ---
$ws.Range("F3").Value2 = "Cronograma do TCC"

# --- 3. Fill in the missing "X" marks on table 1 ---
$ws.Cells.Item(5, 8).Value2 = "X"   # H5
$ws.Cells.Item(7, 8).Value2 = "X"   # H7
$ws.Cells.Item(8, 8).Value2 = "X"   # H8
$ws.Cells.Item(10, 9).Value2 = "X"  # I10
$ws.Cells.Item(11, 9).Value2 = "X"  # I11
$ws.Cells.Item(12, 10).Value2 = "X" # J12

# Wrap text across the whole X-grid of table 1 (matches the source edit).
$ws.Range("G5:J12").WrapText = $true

# --- 4. Re-label the new (second) table's title ---
$ws.Range("F15").Value2 = "Cronograma da Implementação"

# --- 5. Re-label the new table's task names (column F) ---
$ws.Cells.Item(17, 6).Value2 = "1. Implementação de Hardware (Inclui criação do circuito de comunicação Serial)"
$ws.Cells.Item(18, 6).Value2 = "2. Implementação do código no ESP32"
$ws.Cells.Item(19, 6).Value2 = "3. Teste da implementação do código"
$ws.Cells.Item(20, 6).Value2 = "4. Implementação do Aplicativo"
$ws.Cells.Item(21, 6).Value2 = "5. Teste do aplicativo"
# row 22 ("6. Análise dos Resultados") is unchanged from the copy.
$ws.Cells.Item(23, 6).Value2 = "7. Verificar se o projeto atende aos requisitos"
$ws.Cells.Item(24, 6).Value2 = "8. Revisar projeto caso não atenda os requisitos"

# --- 6. Fill in the new table's "X" marks that differ from the copied pattern ---
$ws.Cells.Item(19, 8).Value2 = "X"  # H19
$ws.Cells.Item(20, 8).Value2 = "X"  # H20
$ws.Cells.Item(24, 10).Value2 = "X" # J24

# --- 7. Make sure the new title row is merged (Copy should already do this) ---
if ($ws.Range("F15:J15").MergeCells -ne $true) {
    $ws.Range("F15:J15").Merge()
}

# --- 8. Widen column F to fit the longer task descriptions ---
$ws.Columns.Item(6).ColumnWidth = 63.4

# --- 9. Update the active selection to the new table's title row ---
$ws.Range("F15:J15").Select()
